# Update attachment_filename column (B) values to contain the full absolute
# path to the file, instead of just the bare filename.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$basePath = "D:\Users\Chickens\Documents\EPCC\SynthSys\code_projects\synbio-toolkit\src\test\resources\ed\biordm\sbol\synbio\handler\"

$ws.Range("B2").Value = $basePath + "NC_001499.gbk"
$ws.Range("B3").Value = $basePath + "NC_035470.gbk"
$ws.Range("B5").Value = $basePath + "NC_014139.gbk"

# Move the active selection, mirroring the cursor position left behind in the
# saved file after the edit.
$ws.Range("D16").Select()
